$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 160, shifting existing rows 160-265 down to 161-266.
$ws.Rows("160:160").Insert()

# Populate the newly inserted row 160 with the new data record.
$ws.Range("A160").Value = 7
$ws.Range("B160").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C160").Value = "Ñuble"
$ws.Range("D160").Value = 44673
$ws.Range("E160").Value = 16
$ws.Range("F160").Value = 100112008
$ws.Range("G160").Value = "Coliflor"
$ws.Range("H160").Value = "Sin especificar"
$ws.Range("I160").Value = "Primera"
$ws.Range("J160").Value = 160
$ws.Range("K160").Value = 1100
$ws.Range("L160").Value = 1200
$ws.Range("M160").Value = 1150
$ws.Range("N160").Value = "$/unidad"
$ws.Range("O160").Value = "Provincia de Diguillín"
$ws.Range("P160").Value = 1150
$ws.Range("Q160").Value = 1
$ws.Range("R160").Value = "Hortaliza"
